$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z7").Formula = "=T7/SQRT(3)"
$ws.Range("AA7").Formula = "=U7/SQRT(3)"
$ws.Range("AB7").Formula = "=V7/SQRT(3)"
$ws.Range("AC7").Formula = "=W7/SQRT(3)"
$ws.Range("AD7").Formula = "=X7/SQRT(3)"

Write-Host "Z7:" $ws.Range("Z7").Formula
Write-Host "Z7 value:" $ws.Range("Z7").Value
